$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2504.3641872458752, 1680.6863481769783, 1660.8938869597118),
    @(2262.55694124916, 1504.6477254743138, 1375.3170900420705),
    @(2560.3680528721502, 1740.2667114603717, 1585.0832652976844),
    @(2494.2847112728741, 1887.6398329156857, 1890.9331066616555),
    @(2528.8038454788189, 1753.003727639835, 1754.5421504869457),
    @(2478.7459565060203, 1828.9944935750555, 1882.8313575376681),
    @(2397.2095768963327, 1849.4812563282758, 1679.6191375567414),
    @(2467.9970013459274, 1930.0372423531173, 1785.1941432769693),
    @(2643.8788154272052, 1943.475416109985, 1635.4938427291459),
    @(2394.9100913362677, 1530.383635317052, 1428.9671136278478),
    @(2163.0424132805206, 1575.4338612267611, 1407.2526586545646),
    @(2770.6264489143277, 2193.8783156302029, 1875.1465558948457),
    @(2528.0865774114127, 1958.223367122167, 1765.5967950965151),
    @(2609.7137253488277, 2022.9225150705993, 1756.5472354606577),
    @(2498.9596433937973, 2031.0777341711305, 1838.001306978955),
    @(2592.3756656780934, 1796.2257130653136, 1555.3439003568135),
    @(2355.9379710138674, 1789.6481733224953, 1700.6493438790021),
    @(2615.4880450976702, 2163.4214778065211, 2077.8516997418374),
    @(2011.2299318360142, 2033.9311347760126, 1981.0688096723513),
    @(2576.2750779554822, 1984.8956900048479, 1877.125556341239),
    @(2713.0795141164522, 2022.0638328393798, 1981.0244367991929),
    @(2569.5114444116962, 1982.4537065438044, 1731.7631208362495),
    @(2527.0568162501113, 1686.8929203187608, 1775.9525937170963),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$ws.Range("A1:C23").Select()
